# Auto-generated edit script applying updated transition-matrix probabilities
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2145454545454545
$ws.Range("C2").Value = 0.5272727272727272
$ws.Range("J2").Value = 0.01090909090909091
$ws.Range("P2").Value = 0.1236363636363636
$ws.Range("S2").Value = 0.1236363636363636
$ws.Range("B3").Value = 0.01282051282051282
$ws.Range("C3").Value = 0.04487179487179487
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.782051282051282
$ws.Range("S3").Value = 0.1346153846153846
$ws.Range("J4").Value = 0.0851063829787234
$ws.Range("P4").Value = 0.574468085106383
$ws.Range("S4").Value = 0.3404255319148936
$ws.Range("B6").Value = 0.08196721311475409
$ws.Range("D6").Value = 0.01229508196721311
$ws.Range("F6").Value = 0.08196721311475409
$ws.Range("J6").Value = 0.1926229508196721
$ws.Range("O6").Value = 0.02049180327868852
$ws.Range("Q6").Value = 0.1270491803278689
$ws.Range("R6").Value = 0.08196721311475409
$ws.Range("S6").Value = 0.4016393442622951
$ws.Range("B7").Value = 0.09722222222222222
$ws.Range("D7").Value = 0.02314814814814815
$ws.Range("F7").Value = 0.06481481481481481
$ws.Range("J7").Value = 0.1435185185185185
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.1898148148148148
$ws.Range("R7").Value = 0.1018518518518518
$ws.Range("S7").Value = 0.3611111111111111
$ws.Range("B8").Value = 0.08531746031746032
$ws.Range("D8").Value = 0.02182539682539682
$ws.Range("E8").Value = 0.001984126984126984
$ws.Range("F8").Value = 0.05158730158730158
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.01785714285714286
$ws.Range("Q8").Value = 0.1646825396825397
$ws.Range("R8").Value = 0.08531746031746032
$ws.Range("S8").Value = 0.4464285714285715
$ws.Range("B9").Value = 0.07627118644067797
$ws.Range("D9").Value = 0.01694915254237288
$ws.Range("F9").Value = 0.05084745762711865
$ws.Range("J9").Value = 0.1228813559322034
$ws.Range("O9").Value = 0.01694915254237288
$ws.Range("Q9").Value = 0.173728813559322
$ws.Range("R9").Value = 0.1144067796610169
$ws.Range("S9").Value = 0.4279661016949153
$ws.Range("B10").Value = 0.09377526273241714
$ws.Range("D10").Value = 0.02021018593371059
$ws.Range("E10").Value = 0.0008084074373484236
$ws.Range("F10").Value = 0.08407437348423606
$ws.Range("J10").Value = 0.1156022635408246
$ws.Range("O10").Value = 0.01535974130962005
$ws.Range("Q10").Value = 0.2101859337105901
$ws.Range("R10").Value = 0.07518189167340339
$ws.Range("S10").Value = 0.3848019401778496
$ws.Range("F11").Value = 0.003546099290780142
$ws.Range("G11").Value = 0.1028368794326241
$ws.Range("J11").Value = 0.07446808510638298
$ws.Range("K11").Value = 0.1631205673758865
$ws.Range("L11").Value = 0.6418439716312057
$ws.Range("S11").Value = 0.01418439716312057
$ws.Range("G12").Value = 0.7853403141361257
$ws.Range("J12").Value = 0.1308900523560209
$ws.Range("K12").Value = 0.03141361256544502
$ws.Range("L12").Value = 0.03141361256544502
$ws.Range("S12").Value = 0.02094240837696335
$ws.Range("G13").Value = 0.7368421052631579
$ws.Range("J13").Value = 0.1929824561403509
$ws.Range("S13").Value = 0.07017543859649122
$ws.Range("F15").Value = 0.02966101694915254
$ws.Range("H15").Value = 0.1186440677966102
$ws.Range("I15").Value = 0.05508474576271186
$ws.Range("J15").Value = 0.3559322033898305
$ws.Range("K15").Value = 0.08050847457627118
$ws.Range("M15").Value = 0.01271186440677966
$ws.Range("O15").Value = 0.1016949152542373
$ws.Range("S15").Value = 0.2457627118644068
$ws.Range("F16").Value = 0.02209944751381215
$ws.Range("H16").Value = 0.1657458563535912
$ws.Range("I16").Value = 0.1270718232044199
$ws.Range("J16").Value = 0.3314917127071823
$ws.Range("K16").Value = 0.0718232044198895
$ws.Range("M16").Value = 0.03867403314917127
$ws.Range("O16").Value = 0.08287292817679558
$ws.Range("S16").Value = 0.1602209944751381
$ws.Range("F17").Value = 0.013215859030837
$ws.Range("H17").Value = 0.1607929515418502
$ws.Range("I17").Value = 0.1299559471365639
$ws.Range("J17").Value = 0.3854625550660793
$ws.Range("K17").Value = 0.08590308370044053
$ws.Range("M17").Value = 0.01101321585903084
$ws.Range("O17").Value = 0.05947136563876652
$ws.Range("S17").Value = 0.1541850220264317
$ws.Range("F18").Value = 0.01477832512315271
$ws.Range("H18").Value = 0.2068965517241379
$ws.Range("I18").Value = 0.1231527093596059
$ws.Range("J18").Value = 0.3251231527093596
$ws.Range("K18").Value = 0.1182266009852217
$ws.Range("M18").Value = 0.03448275862068965
$ws.Range("O18").Value = 0.0541871921182266
$ws.Range("S18").Value = 0.1231527093596059
$ws.Range("F19").Value = 0.01325889741800419
$ws.Range("H19").Value = 0.2309839497557571
$ws.Range("I19").Value = 0.08164689462665736
$ws.Range("J19").Value = 0.3489183531053733
$ws.Range("K19").Value = 0.09211444521981857
$ws.Range("M19").Value = 0.02581995812979763
$ws.Range("N19").Value = 0.001395673412421493
$ws.Range("O19").Value = 0.0628053035589672
$ws.Range("S19").Value = 0.1430565247732031
